$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.186.48'
$ws.Range("E2").Value = '  -1.17%  '
$ws.Range("D3").Value = '1.573.48'
$ws.Range("E3").Value = '  -0.48%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = "'207.58"
$ws.Range("E5").Value = '  -0.02%  '
$ws.Range("D6").Value = "'0.490"
$ws.Range("E6").Value = '  -1.71%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = "'22.26"
$ws.Range("E8").Value = '  +0.20%  '
$ws.Range("E9").Value = '  -0.73%  '
$ws.Range("E10").Value = '  +0.11%  '
$ws.Range("E11").Value = '  +0.20%  '
$ws.Range("D12").Value = '1.796.33'
$ws.Range("E12").Value = '  -0.45%  '
$ws.Range("D13").Value = '1.574.60'
$ws.Range("E13").Value = '  -0.52%  '
$ws.Range("D14").Value = "'3.77"
$ws.Range("E14").Value = '  -1.25%  '
$ws.Range("D15").Value = "'0.519"
$ws.Range("E15").Value = '  -0.95%  '
$ws.Range("D16").Value = '27.192.92'
$ws.Range("E16").Value = '  -1.10%  '
$ws.Range("D17").Value = "'62.16"
$ws.Range("E18").Value = '  +0.95%  '
$ws.Range("D19").Value = "'214.07"
$ws.Range("E19").Value = '  -0.50%  '
$ws.Range("E20").Value = '  -0.69%  '
$ws.Range("E21").Value = '  -0.05%  '
$ws.Range("E22").Value = '  -0.18%  '
$ws.Range("D23").Value = "'9.39"
$ws.Range("E23").Value = '  -3.29%  '
$ws.Range("E24").Value = '  +0.28%  '
$ws.Range("E25").Value = '  -0.26%  '
$ws.Range("D26").Value = "'6.68"
$ws.Range("E26").Value = '  -3.57%  '
$ws.Range("D27").Value = "'14.95"
$ws.Range("E27").Value = '  -0.61%  '
$ws.Range("E28").Value = '  -0.03%  '
$ws.Range("E29").Value = '  -1.25%  '
$ws.Range("D30").Value = "'1.12"
$ws.Range("E30").Value = '  -2.62%  '
$ws.Range("E31").Value = '  -1.63%  '
$ws.Range("E32").Value = '  -1.16%  '
$ws.Range("D33").Value = '1.409.18'
$ws.Range("E33").Value = '  +3.23%  '
$ws.Range("E34").Value = '  -1.28%  '
$ws.Range("E35").Value = '  +1.82%  '
$ws.Range("E36").Value = '  -0.97%  '
$ws.Range("D37").Value = "'0.937"
$ws.Range("E37").Value = '  -3.44%  '
$ws.Range("E38").Value = '  -1.73%  '
$ws.Range("D39").Value = "'0.818"
$ws.Range("E39").Value = '  -0.50%  '
$ws.Range("E40").Value = '  -2.59%  '
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("D42").Value = "'0.994"
$ws.Range("E42").Value = '  +2.41%  '
$ws.Range("D43").Value = "'1.83"
$ws.Range("E43").Value = '  +4.13%  '
$ws.Range("D44").Value = "'5.38"
$ws.Range("E44").Value = '  +2.27%  '
$ws.Range("E45").Value = '  +0.93%  '
$ws.Range("D46").Value = "'63.74"
$ws.Range("E46").Value = '  -0.58%  '
$ws.Range("D47").Value = '1.709.39'
$ws.Range("E47").Value = '  -0.42%  '
$ws.Range("D48").Value = "'85.98"
$ws.Range("E48").Value = '  -0.46%  '
$ws.Range("D49").Value = '0.0₇0993'
$ws.Range("E49").Value = '  -0.39%  '
$ws.Range("E50").Value = '  -0.38%  '
$ws.Range("D51").Value = "'0.0493"
$ws.Range("E51").Value = '  -0.06%  '

# Reset style on quote-prefixed text cells so no stray quotePrefix/style is visibly applied
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D51").Style = "Normal"
